$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @("60+2=62","59-21=38","32+62=94","15+62=77","43+33=76","21+78=99","34+1=35","86-47=39","12+36=48","88-65=23","49+20=69","92+3=95","69-69=0","24+20=44","0+75=75","93+1=94","86+11=97","69-0=69","92-79=13","92-6=86","70-45=25","55-6=49","27-18=9","99-32=67","59-31=28","67-66=1","44+21=65","16+4=20","12+3=15","84-51=33","5+21=26","96-33=63","10+13=23","56-8=48","83+12=95","44-6=38","37-15=22","77-57=20","38-7=31","3+69=72","38+57=95","43-41=2","27+43=70","97-64=33","15+6=21","61-30=31","49+6=55","52+15=67","21+20=41","38+32=70","82-29=53","88+10=98","58-29=29","78-70=8","7+78=85","92-75=17","54+9=63","19+56=75","2+32=34","91-57=34","80-47=33","60+34=94","71-33=38","38+26=64","37+16=53","85+3=88","91-29=62","10-3=7","43-41=2","50+26=76","10+85=95","71+27=98","58-25=33","29+45=74","18+65=83","99-43=56","22+26=48","55+28=83","4+15=19","13+50=63","50-34=16","14-10=4","77-58=19","63-44=19","54-43=11","33+2=35","99-82=17","80-36=44","74+1=75","18-2=16","28-6=22","85-46=39","85+11=96","99-56=43","15-8=7","35+20=55","41+3=44","32+23=55","52+5=57","93+3=96")
$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $rng = $cell.Range
        $rng.End = $rng.End - 1
        $rng.Text = $values[$idx]
        $idx = $idx + 1
    }
}
Write-Host "Done. idx=" $idx
